$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 232 (pushes existing row 232.. down to 233..)
$ws.Rows(232).Insert()

# Populate the newly inserted row 232 with the new weekly record
$ws.Cells.Item(232, 1).Value = 3
$ws.Cells.Item(232, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(232, 3).Value = "Coquimbo"
$ws.Cells.Item(232, 4).Value = 44523
$ws.Cells.Item(232, 5).Value = 5
$ws.Cells.Item(232, 6).Value = 100112037
$ws.Cells.Item(232, 7).Value = "Cebollín"
$ws.Cells.Item(232, 8).Value = "Sin especificar"
$ws.Cells.Item(232, 9).Value = "Primera"
$ws.Cells.Item(232, 10).Value = 160
$ws.Cells.Item(232, 11).Value = 3000
$ws.Cells.Item(232, 12).Value = 3000
$ws.Cells.Item(232, 13).Value = 3000
$ws.Cells.Item(232, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(232, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(232, 16).Value = 83
$ws.Cells.Item(232, 17).Value = 36
$ws.Cells.Item(232, 18).Value = "Hortaliza"
